$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its original location (end of the
#    "Project Title" paragraph). It gets re-added later, at the very end
#    of the document.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 2. Merge the four runs of the second bulleted question into a single
#    run with the concatenated text.
# ---------------------------------------------------------------------
$ldq = [char]0x201C
$rdq = [char]0x201D
$questionText = "Does general user base sentiment on twitter about certain companies correlate to overall company performance as measured through their stock price?  Is this better or worse than the " + $ldq + "experts" + $rdq + "?"
$d.Content.Find.Execute($questionText, $true, $false, $false, $false, $false, $true, 1, $false, $questionText, 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Replace the last two blank trailing paragraphs with the new
#    "BACKUP QUESTIONS TO ANSWER" section.
# ---------------------------------------------------------------------

# Paragraph 19 (first of the three trailing blank paragraphs) is left
# untouched. Paragraph 20 becomes the yellow-highlighted heading.
$pHeading = $d.Paragraphs(20)
$pHeading.Range.Text = "BACKUP QUESTIONS TO ANSWER:"
$pHeading.Range.HighlightColorIndex = 7   # wdYellow

# Paragraph 21 becomes the first bulleted item of the new list (ilvl 0).
# Applying a list template from the bullet gallery mints a brand new
# numId/abstractNum pair (bullet glyphs at 3 rotating levels using
# Symbol / Courier New / Wingdings, same shape as a fresh Word bullet
# list).
$bulletTemplate = $word.ListGalleries.Item(1).ListTemplates.Item(1)

$idx = 21
$p = $d.Paragraphs($idx)
$p.Range.Text = "Using historical Kaggle dataset:"
$p.Style = "List Paragraph"
$p.Range.ListFormat.ApplyListTemplate($bulletTemplate)

# --- remaining bullets of the backup list, all sharing numId 5 -------
$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>What is the average amount of time for maximum profit?</w:t></w:r></w:p>'
$d.Paragraphs($idx).Range.InsertXML($xml) | Out-Null

$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Aka is time in the market more profitable than timing the market?</w:t></w:r></w:p>'
$d.Paragraphs($idx).Range.InsertXML($xml) | Out-Null

$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$falsePeaksQ = "How many " + $ldq + "false peaks" + $rdq + " occurred during this maximum profit period?"
$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>' + $falsePeaksQ + '</w:t></w:r></w:p>'
$d.Paragraphs($idx).Range.InsertXML($xml) | Out-Null

$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>False peaks would be where the stock price was going up and then started going down.  Could count the number of new maxes before reaching the final one</w:t></w:r></w:p>'
$d.Paragraphs($idx).Range.InsertXML($xml) | Out-Null

$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Story: if large </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>amount</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> of peaks before final peak then this is more attributed to luck than anything</w:t></w:r></w:p>'
$d.Paragraphs($idx).Range.InsertXML($xml) | Out-Null

$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Correlate to SEC data?</w:t></w:r></w:p>'
$d.Paragraphs($idx).Range.InsertXML($xml) | Out-Null

# Final bullet: carries the restored _GoBack bookmark plus a trailing
# line break inside its own run.
$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Maximum change day over day and does this correlate with SEC data for that stock?</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:br/></w:r></w:p>'
$d.Paragraphs($idx).Range.InsertXML($xml) | Out-Null

Write-Output "done"
